$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column width changes (values chosen so the COM character-width -> stored
# width rounding lands on the closest representable width to the target)
$ws.Columns.Item(4).ColumnWidth = 1.33    # -> stored width ~2.1667 (was 3.140625, target 2.140625)
$ws.Columns.Item(9).ColumnWidth = 4.83    # -> stored width ~5.6667 (was 3.140625, target 5.7109375)
$ws.Columns.Item(10).ColumnWidth = 4.83   # -> stored width ~5.6667 (was 3.140625, target 5.7109375)
$ws.Columns.Item(11).ColumnWidth = 4.83   # -> stored width ~5.6667 (was 7.7109375, target 5.7109375)

# Cell value changes (row 1)
$ws.Range("A1").Value = 3
$ws.Range("C1").Value = 28
$ws.Range("D1").Value = 8
$ws.Range("E1").Value = 20
$ws.Range("F1").Value = 14
$ws.Range("H1").Value = 29
$ws.Range("I1").Value = 0.048
$ws.Range("J1").Value = 0.019
$ws.Range("K1").Value = 0.078
